# Applies:
#  - remove the pStyle/rStyle formatting (teibibl / teititlem) from the
#    lone body paragraph, leaving just the bookmark
#  - drop the duplicate "egXML"/"egXMLTable" style pair and rename the
#    surviving "egXML0"/"egXMLTable0" pair back to "egXML"/"egXMLTable"
#  - add the new "marginRight" / "marginLeft" paragraph styles (based on
#    marginOuter) to support @cRef on <ptr>

$d = $word.ActiveDocument

# --- document.xml: strip the paragraph formatting pointing at teibibl/teititlem ---
$p = $d.Paragraphs(1)
$p.Style = "Normal"

# --- styles.xml: collapse egXML/egXMLTable duplication ---
# The template carries two copies of these custom styles: the original
# "egXML"/"egXMLTable" pair, and a second pair that got minted as
# "egXML0"/"egXMLTable0" (same w:name, different w:styleId) to avoid a
# collision. Drop the original pair, then delete+recreate the "0" pair so
# it lands on the clean "egXML"/"egXMLTable" ids.
$d.Styles("egXML").Delete()
$d.Styles("egXMLTable").Delete()

$d.Styles("egXML0").Delete()
$egXML = $d.Styles.Add("egXML", 1)
$egXML.BaseStyle = "Normal"
$egXML.QuickStyle = $true
$egXML.Font.Name = "Courier"
$egXML.Font.Size = 10

$d.Styles("egXMLTable0").Delete()
$egXMLTable = $d.Styles.Add("egXMLTable", 1)
$egXMLTable.BaseStyle = "Normal"
$egXMLTable.QuickStyle = $true
$egXMLTable.ParagraphFormat.SpaceBefore = 4
$egXMLTable.Font.Name = "Courier"
$egXMLTable.Font.Size = 9

# --- styles.xml: add marginRight / marginLeft (based on marginOuter) ---
$marginRight = $d.Styles.Add("marginRight", 1)
$marginRight.BaseStyle = "marginOuter"
$marginRight.NextParagraphStyle = "Normal"
$marginRight.QuickStyle = $true

$marginLeft = $d.Styles.Add("marginLeft", 1)
$marginLeft.BaseStyle = "marginOuter"
$marginLeft.NextParagraphStyle = "Normal"
$marginLeft.QuickStyle = $true
